$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.402482867240906
$ws.Range("B1").Value = 1.616135835647583
$ws.Range("C1").Value = 5.997123718261719
$ws.Range("D1").Value = 2.654081344604492
$ws.Range("E1").Value = 1.166913986206055
